$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 39/40 swap: Filecoin <-> EnergySwap ---
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "26.56"
$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "4.61"
$ws.Range("E40").Value = "  +0.05%  "

# --- Simple value updates ---
$ws.Range("D2").Value = "68.105.89"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "3.248.19"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "582.20"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").Value = "184.63"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "0.601"
$ws.Range("E8").Value = "  +0.73%  "
$ws.Range("E9").Value = "  -3.23%  "
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "3.813.31"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("D14").Value = "27.94"
$ws.Range("E14").Value = "  -2.84%  "
$ws.Range("D15").Value = "68.129.18"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "0.0000170"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").Value = "3.220.27"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").Value = "13.49"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("D20").Value = "396.21"
$ws.Range("E20").Value = "  +4.33%  "
$ws.Range("D21").Value = "7.61"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").Value = "71.40"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("D26").Value = "0.187"
$ws.Range("E26").Value = "  +2.82%  "
$ws.Range("D27").Value = "9.65"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("D30").Value = "5.62"
$ws.Range("E30").Value = "  -1.31%  "
$ws.Range("D31").Value = "22.83"
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("D32").Value = "7.04"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("E35").Value = "  -5.17%  "
$ws.Range("D36").Value = "162.11"
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("E37").Value = "  +2.09%  "
$ws.Range("D38").Value = "0.816"
$ws.Range("E38").Value = "  -3.25%  "
$ws.Range("D41").Value = "6.50"
$ws.Range("E41").Value = "  -2.53%  "
$ws.Range("E42").Value = "  -4.45%  "
$ws.Range("D43").Value = "41.24"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("E44").Value = "  -0.66%  "
$ws.Range("D45").Value = "25.05"
$ws.Range("E45").Value = "  -1.99%  "
$ws.Range("D46").Value = "2.611.82"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("D47").Value = "335.54"
$ws.Range("E47").Value = "  -3.59%  "
$ws.Range("E48").Value = "  -2.44%  "
$ws.Range("D49").Value = "6.34"
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("E51").Value = "  -1.20%  "
